$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "JDE"
$ws.Name = "JDE"

# Move the data block from F1:H31 to C1:E31 (copy then clear the source)
$src = $ws.Range("F1:H31")
$dst = $ws.Range("C1:E31")
$src.Copy($dst) | Out-Null
$src.Clear() | Out-Null

# Restore the view: top-left at A1, selection on the header row of the
# relocated table (C1:E1), matching where the data now lives.
$ws.Range("C1:E1").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
